{"js": "// Load all paragraphs in the document body.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Helper to find the (first) paragraph whose text starts with the given\n// prefix, so we stay robust to the exact paragraph index.\nfunction findParagraph(items, prefix) {\n  for (let i = 0; i < items.length; i++) {\n    if (items[i].text.indexOf(prefix) === 0) {\n      return items[i];\n    }\n  }\n  return null;\n}\n\n// 1) Update the activation date: 01/01/2012 -> 01/01/2023\nconst dateResults = body.search(\"Ativa\u00e7\u00e3o: 01/01/2012\", { matchCase: true });\ndateResults.load(\"text\");\nawait context.sync();\nif (dateResults.items.length > 0) {\n  dateResults.items[0].insertText(\"Ativa\u00e7\u00e3o: 01/01/2023\", \"Replace\");\n}\n\n// 2) Insert the English translation of \"Objetivos\" right after the\n//    Portuguese paragraph.\nconst objetivosPt = findParagraph(\n  paragraphs.items,\n  \"Apresentar uma vis\u00e3o geral da qu\u00edmica dos elementos\"\n);\nif (objetivosPt) {\n  const objetivosEn = objetivosPt.insertParagraph(\n    \"To present an overview of the chemistry of elements and their compounds, emphasizing the correlations between physical and chemical properties with structural and binding aspects, methods of obtaining them in laboratory and industry, in addition to the main properties and applications.\",\n    \"After\"\n  );\n  objetivosEn.font.italic = true;\n}\n\n// 3) Insert the English translation of \"Programa resumido\" right after the\n//    Portuguese paragraph (the short summary).\nconst resumidoPt = findParagraph(\n  paragraphs.items,\n  \"Ocorr\u00eancia, obten\u00e7\u00e3o, estrutura, propriedades e aplica\u00e7\u00f5es de elementos met\u00e1licos e n\u00e3o-met\u00e1licos\"\n);\nif (resumidoPt) {\n  const resumidoEn = resumidoPt.insertParagraph(\n    \"Occurrence, obtaining, structure, properties and applications of metallic and non-metallic elements; polyatomic molecules; halogenated compounds and the oxygen, nitrogen, carbon and boron families; oxygenated compounds. Industrial manufacturing processes.\",\n    \"After\"\n  );\n  resumidoEn.font.italic = true;\n}\n\n// 4) Insert the English translation of \"Programa\" right after the\n//    Portuguese paragraph (the full program).\nconst programaPt = findParagraph(\n  paragraphs.items,\n  \"Ocorr\u00eancia, obten\u00e7\u00e3o, estrutura, propriedades e aplica\u00e7\u00f5es de elementos n\u00e3o-met\u00e1licos\"\n);\nif (programaPt) {\n  const programaEn = programaPt.insertParagraph(\n    \"Occurrence, obtaining, structure, properties and applications of non-metallic elements: noble gases, molecular hydrogen, halogens, molecular oxygen, ozone and molecular nitrogen; semimetals; alkali, alkaline earth and transition metals; polyatomic molecules and catenated species of: sulfur, phosphorus and carbon; halogenated compounds and the oxygen, nitrogen, carbon and boron families; oxygenated compounds: oxides and oxycompounds. Industrial manufacturing processes of the main chemical inputs and materials.\",\n    \"After\"\n  );\n  programaEn.font.italic = true;\n}\n\n// 5) Update the prerequisite line.\nconst reqResults = body.search(\"LOQ4031 -  Qu\u00edmica Geral I  (Requisito)\", {\n  matchCase: true,\n});\nreqResults.load(\"text\");\nawait context.sync();\nif (reqResults.items.length > 0) {\n  reqResults.items[0].insertText(\n    \"LOQ4100 -  Fundamentos de Qu\u00edmica para Engenharia I (Requisito fraco)\",\n    \"Replace\"\n  );\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Update the activation date: 01/01/2012 -> 01/01/2023\n$find = $d.Content.Find\n$find.Execute(\n    \"Ativa\u00e7\u00e3o: 01/01/2012\",\n    $false, $false, $false, $false, $false, $true, 1, $false,\n    \"Ativa\u00e7\u00e3o: 01/01/2023\",\n    2\n)\n\n# Helper: find the index (1-based, as used by Paragraphs.Item) of the first\n# paragraph whose text starts with the given prefix.\nfunction Get-ParagraphIndex($doc, $prefix) {\n    $i = 0\n    foreach ($p in $doc.Paragraphs) {\n        $i = $i + 1\n        if ($p.Range.Text.StartsWith($prefix)) {\n            return $i\n        }\n    }\n    return -1\n}\n\n# Helper: insert a new italic paragraph right after the paragraph at\n# 1-based index $index, containing $text.\nfunction Add-ItalicParagraphAfter($doc, $index, $text) {\n    $para = $doc.Paragraphs.Item($index)\n    $para.Range.InsertParagraphAfter()\n    $newPara = $doc.Paragraphs.Item($index + 1)\n    $newPara.Range.Text = $text\n    $newRange = $doc.Range($newPara.Range.Start, $newPara.Range.Start + $text.Length)\n    $newRange.Font.Italic = $true\n}\n\n# 2) Insert the English translation of \"Objetivos\" right after the\n#    Portuguese paragraph.\n$objetivosIdx = Get-ParagraphIndex $d \"Apresentar uma vis\u00e3o geral da qu\u00edmica dos elementos\"\nAdd-ItalicParagraphAfter $d $objetivosIdx \"To present an overview of the chemistry of elements and their compounds, emphasizing the correlations between physical and chemical properties with structural and binding aspects, methods of obtaining them in laboratory and industry, in addition to the main properties and applications.\"\n\n# 3) Insert the English translation of \"Programa resumido\" right after the\n#    Portuguese paragraph (the short summary).\n$resumidoIdx = Get-ParagraphIndex $d \"Ocorr\u00eancia, obten\u00e7\u00e3o, estrutura, propriedades e aplica\u00e7\u00f5es de elementos met\u00e1licos e n\u00e3o-met\u00e1licos\"\nAdd-ItalicParagraphAfter $d $resumidoIdx \"Occurrence, obtaining, structure, properties and applications of metallic and non-metallic elements; polyatomic molecules; halogenated compounds and the oxygen, nitrogen, carbon and boron families; oxygenated compounds. Industrial manufacturing processes.\"\n\n# 4) Insert the English translation of \"Programa\" right after the\n#    Portuguese paragraph (the full program).\n$programaIdx = Get-ParagraphIndex $d \"Ocorr\u00eancia, obten\u00e7\u00e3o, estrutura, propriedades e aplica\u00e7\u00f5es de elementos n\u00e3o-met\u00e1licos\"\nAdd-ItalicParagraphAfter $d $programaIdx \"Occurrence, obtaining, structure, properties and applications of non-metallic elements: noble gases, molecular hydrogen, halogens, molecular oxygen, ozone and molecular nitrogen; semimetals; alkali, alkaline earth and transition metals; polyatomic molecules and catenated species of: sulfur, phosphorus and carbon; halogenated compounds and the oxygen, nitrogen, carbon and boron families; oxygenated compounds: oxides and oxycompounds. Industrial manufacturing processes of the main chemical inputs and materials.\"\n\n# 5) Update the prerequisite line.\n$find2 = $d.Content.Find\n$find2.Execute(\n    \"LOQ4031 -  Qu\u00edmica Geral I  (Requisito)\",\n    $false, $false, $false, $false, $false, $true, 1, $false,\n    \"LOQ4100 -  Fundamentos de Qu\u00edmica para Engenharia I (Requisito fraco)\",\n    2\n)\n"}
